$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7..44 down to 8..45
$ws.Rows(7).Insert()

# Copy the formatting of the row above (row 6) into the new row 7 so the
# S.No. cell (A7) keeps the same bold/bordered/centered style as the rest
# of column A.
$ws.Range("A6").Copy($ws.Range("A7"))

# Populate the newly inserted row with the new book's data
$ws.Range("A7").Value2 = 6
$ws.Range("B7").Value2 = "An Introduction to Analysis - Gerald G. Bilodeau, Paul R. Thie, G. E. Kenough (2010, Jones and Bartlett) 2nd Edition.pdf"
$ws.Range("C7").Value2 = "An Introduction to Analysis"
$ws.Range("D7").Value2 = "Gerald G. Bilodeau, Paul R. Thie, G. E. Kenough"
$ws.Range("E7").Value2 = "[Drive](https://drive.google.com/file/d/1t3svnhU23xaq6IcewyC_F7qoNYv4eEAB/view)"
$ws.Range("F7").Value2 = "2nd"
$ws.Range("G7").Value2 = 2010
$ws.Range("H7").Value2 = "Jones and Bartlett"

# The S.No. column (A) is a plain sequential counter (1..N) and the Insert
# above left the old numbers in place on rows 8..45 (they did not shift
# along with the row insert). Renumber the whole column so it reads
# 1,2,3,...,44 again.
for ($i = 2; $i -le 45; $i++) {
    $ws.Cells.Item($i, 1).Value2 = $i - 1
}
